# Update the "dSF" column (F) values for specific rows to reflect
# re-pulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    4  = -1
    9  = 0
    11 = 0
    16 = 2
    17 = 3
    22 = 0
    25 = -2
    35 = 0
    36 = 0
    51 = 1
    53 = 0
    55 = 1
    56 = -1
    60 = 0
    62 = 0
    70 = -1
    72 = -2
    74 = 0
    77 = 1
    78 = -1
    81 = 1
    83 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
